$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Workbook-level view state: NOV-2020 becomes the active/selected sheet tab.
# ---------------------------------------------------------------------------
$wsOct = $wb.Worksheets.Item("OCT-2020")
$wsNov = $wb.Worksheets.Item("NOV-2020")

# ---------------------------------------------------------------------------
# OCT-2020 ("sheet2"): tab was selected before, now it's not; the old
# selection/active cell is cleared back to a plain A1:J4 selection.
# ---------------------------------------------------------------------------
$wsOct.Select()
$wsOct.Range("A1:J4").Select()

# ---------------------------------------------------------------------------
# NOV-2020 ("sheet3"): this is the sheet that actually received new task
# rows + a full re-alignment (center -> left) of the tracker table.
# ---------------------------------------------------------------------------
$wsNov.Select()

# New task entries for 2-Nov-2020.
$wsNov.Range("A2").Value = 1
$wsNov.Range("B2").Value = [double]44137
$wsNov.Range("C2").Value = "RPA SONY"
$wsNov.Range("D2").Value = "Conversion issue during CSV files creation at Sony Daily SchedulingReport task has been fixed, tested and running smoothly for all 11 files"
$wsNov.Range("E2").Value = 1
$wsNov.Range("F2").Value = "Completed"

$wsNov.Range("A3").Value = 2
$wsNov.Range("B3").Value = [double]44137
$wsNov.Range("C3").Value = "RPA SAMSUNG"
$wsNov.Range("D3").Value = "Upload issue at RETURNCREDIT MONTHLY task has been fixed, tested and running smoothly"
$wsNov.Range("E3").Value = 1
$wsNov.Range("F3").Value = "Completed"

$wsNov.Range("A4").Value = 3
$wsNov.Range("B4").Value = [double]44137
$wsNov.Range("C4").Value = "RPA SAMSUNG"
$wsNov.Range("D4").Value = "Conversion issue during CSV files creation at SAWDISCOUNT Monthly task is going on"
$wsNov.Range("E4").Value = 0.8
$wsNov.Range("F4").Value = "WIP"

# Whole tracker table (header + data + blank rows through 16) switches from
# center-aligned to left-aligned.
$wsNov.Range("A1:G16").HorizontalAlignment = -4131

# Legend block (rows 19-23) also switches to left alignment.
$wsNov.Range("B19:C23").HorizontalAlignment = -4131

# Row heights for the 3 newly-populated rows match the rest of the table.
$wsNov.Rows(3).RowHeight = 26.4
$wsNov.Rows(4).RowHeight = 28.2

# Column widths / best-fit sizing to match the new content.
$wsNov.Columns("A").ColumnWidth = 3.44140625
$wsNov.Columns("B").ColumnWidth = 10.33203125
$wsNov.Columns("C").ColumnWidth = 13.21875
$wsNov.Columns("D").ColumnWidth = 77.21875
$wsNov.Columns("E").ColumnWidth = 14.77734375
$wsNov.Columns("F").ColumnWidth = 9.88671875
$wsNov.Columns("G").ColumnWidth = 10.109375

$wsNov.Range("D7").Select()

Write-Host "done"
